# Update "想去人数" (want-to-go count) figures in the 苏州-漫展信息 workbook
# to match newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 16003
$ws1.Range("F9").Value = 15517
$ws1.Range("F18").Value = 213
$ws1.Range("F28").Value = 27
$ws1.Range("F39").Value = 5621

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 16003
$ws4.Range("F9").Value = 15517
$ws4.Range("F18").Value = 213
$ws4.Range("F28").Value = 27
$ws4.Range("F41").Value = 5621
